$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.049.21"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").Value = "3.822.29"
$ws.Range("E3").Value = "  +8.00%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'428.81"
$ws.Range("E5").Value = "  +8.25%  "
$ws.Range("D6").Value = "'131.42"
$ws.Range("E6").Value = "  +4.17%  "
$ws.Range("D7").Value = "3.818.44"
$ws.Range("E7").Value = "  +8.13%  "
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'0.735"
$ws.Range("E10").Value = "  +6.93%  "
$ws.Range("D12").Value = "'0.0000339"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "'41.65"
$ws.Range("E13").Value = "  +5.66%  "
$ws.Range("D14").Value = "'10.52"
$ws.Range("E14").Value = "  +13.11%  "
$ws.Range("D15").Value = "4.434.38"
$ws.Range("E15").Value = "  +8.51%  "
$ws.Range("D16").Value = "'15.50"
$ws.Range("E16").Value = "  +21.65%  "
$ws.Range("D18").Value = "3.806.11"
$ws.Range("E18").Value = "  +7.90%  "
$ws.Range("D19").Value = "'20.09"
$ws.Range("E19").Value = "  +6.44%  "
$ws.Range("E20").Value = "  +8.01%  "
$ws.Range("D21").Value = "66.299.35"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").Value = "'416.93"
$ws.Range("E22").Value = "  +4.08%  "
$ws.Range("D23").Value = "'15.19"
$ws.Range("E23").Value = "  +8.42%  "
$ws.Range("D24").Value = "'85.50"
$ws.Range("E24").Value = "  +4.47%  "
$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  +8.29%  "
$ws.Range("D26").Value = "'37.27"
$ws.Range("E26").Value = "  +9.09%  "
$ws.Range("D27").Value = "'10.13"
$ws.Range("E27").Value = "  +14.00%  "
$ws.Range("D28").Value = "'3.31"
$ws.Range("E28").Value = "  +9.68%  "
$ws.Range("D29").Value = "'9.46"
$ws.Range("E29").Value = "  +37.50%  "
$ws.Range("D30").Value = "'5.38"
$ws.Range("E30").Value = "  +2.72%  "
$ws.Range("D31").Value = "'14.29"
$ws.Range("E31").Value = "  +19.68%  "
$ws.Range("D32").Value = "'710.97"
$ws.Range("E32").Value = "  +4.58%  "
$ws.Range("D33").Value = "'0.126"
$ws.Range("E33").Value = "  +12.97%  "
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("D35").Value = "'5.86"
$ws.Range("E35").Value = "  +41.41%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'39.00"
$ws.Range("E37").Value = "  +5.30%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "'55.80"
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").Value = "'0.0472"
$ws.Range("E40").Value = "  +7.24%  "
$ws.Range("D41").Value = "0.0₃0727"
$ws.Range("E41").Value = "  +15.51%  "
$ws.Range("D42").Value = "'2.90"
$ws.Range("E42").Value = "  +2.89%  "
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +9.74%  "
$ws.Range("D46").Value = "'3.24"
$ws.Range("E46").Value = "  +4.97%  "
$ws.Range("D47").Value = "'0.324"
$ws.Range("E47").Value = "  +16.78%  "
$ws.Range("D48").Value = "'2.45"
$ws.Range("E48").Value = "  +43.20%  "
$ws.Range("E49").Value = "  +5.96%  "
$ws.Range("E50").Value = "  +5.48%  "
$ws.Range("D51").Value = "'2.85"
$ws.Range("E51").Value = "  +4.63%  "
